# Scheduled-runner style update: refresh cached market-board figures
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) for the
# Leve profitability tables across all eight job sheets, per the latest
# pull. Values only -- no structural/formatting changes.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 479.8
$ws.Range("I2").Value = 479.8
$ws.Range("K2").Value = 479.8
$ws.Range("M2").Value = -366.8
$ws.Range("H58").Value = 253
$ws.Range("I58").Value = 253
$ws.Range("K58").Value = 759
$ws.Range("M58").Value = -609
$ws.Range("H74").Value = 14139.737
$ws.Range("I74").Value = 14979.706
$ws.Range("K74").Value = 14979.706
$ws.Range("M74").Value = -14043.706
$ws.Range("H77").Value = 14139.737
$ws.Range("I77").Value = 14979.706
$ws.Range("K77").Value = 74898.53
$ws.Range("M77").Value = -70218.53
$ws.Range("H86").Value = 3657.5557
$ws.Range("I86").Value = 2864.75
$ws.Range("K86").Value = 2864.75
$ws.Range("M86").Value = -1741.75
$ws.Range("H87").Value = 160833.17
$ws.Range("J87").Value = 160833.17
$ws.Range("L87").Value = 160833.17
$ws.Range("N87").Value = -163329.17
$ws.Range("H89").Value = 3657.5557
$ws.Range("I89").Value = 2864.75
$ws.Range("K89").Value = 14323.75
$ws.Range("M89").Value = -8707.75
$ws.Range("H90").Value = 160833.17
$ws.Range("J90").Value = 160833.17
$ws.Range("L90").Value = 482499.51
$ws.Range("N90").Value = -494979.51
$ws.Range("H111").Value = 1821.375
$ws.Range("I111").Value = 1293
$ws.Range("J111").Value = 2138.4
$ws.Range("K111").Value = 3879
$ws.Range("L111").Value = 6415.200000000001
$ws.Range("M111").Value = -812
$ws.Range("N111").Value = -12549.2
$ws.Range("H116").Value = 9867.25
$ws.Range("I116").Value = 2984.5
$ws.Range("J116").Value = 16750
$ws.Range("K116").Value = 2984.5
$ws.Range("L116").Value = 16750
$ws.Range("M116").Value = 457.5
$ws.Range("N116").Value = -23634
$ws.Range("H137").Value = 6140.273
$ws.Range("I137").Value = 2790.5652
$ws.Range("J137").Value = 13844.6
$ws.Range("K137").Value = 8371.695599999999
$ws.Range("L137").Value = 41533.8
$ws.Range("M137").Value = -5821.695599999999
$ws.Range("N137").Value = -46633.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1909.3658
$ws.Range("I32").Value = 1932.1
$ws.Range("K32").Value = 1932.1
$ws.Range("M32").Value = -1645.1
$ws.Range("H61").Value = 2372.5405
$ws.Range("I61").Value = 1639.1364
$ws.Range("J61").Value = 3448.2
$ws.Range("K61").Value = 1639.1364
$ws.Range("L61").Value = 3448.2
$ws.Range("M61").Value = -1427.1364
$ws.Range("N61").Value = -3872.2
$ws.Range("H74").Value = 144241.83
$ws.Range("I74").Value = 214954.77
$ws.Range("K74").Value = 214954.77
$ws.Range("M74").Value = -214080.77
$ws.Range("H77").Value = 144241.83
$ws.Range("I77").Value = 214954.77
$ws.Range("K77").Value = 1074773.85
$ws.Range("M77").Value = -1070405.85
$ws.Range("H110").Value = 1266.579
$ws.Range("I110").Value = 1184.6
$ws.Range("J110").Value = 1574
$ws.Range("K110").Value = 1184.6
$ws.Range("L110").Value = 1574
$ws.Range("M110").Value = 860.4000000000001
$ws.Range("N110").Value = -5664
$ws.Range("H132").Value = 2678.3235
$ws.Range("I132").Value = 2463.423
$ws.Range("K132").Value = 7390.268999999999
$ws.Range("M132").Value = -4860.268999999999
$ws.Range("H136").Value = 2372.5405
$ws.Range("I136").Value = 1639.1364
$ws.Range("J136").Value = 3448.2
$ws.Range("K136").Value = 4917.4092
$ws.Range("L136").Value = 10344.6
$ws.Range("M136").Value = -2367.4092
$ws.Range("N136").Value = -15444.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25630
$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27184
$ws.Range("H80").Value = 495.33334
$ws.Range("I80").Value = 329.66666
$ws.Range("K80").Value = 329.66666
$ws.Range("M80").Value = 668.33334
$ws.Range("H83").Value = 495.33334
$ws.Range("I83").Value = 329.66666
$ws.Range("K83").Value = 1648.3333
$ws.Range("M83").Value = 3343.6667
$ws.Range("H86").Value = 1878.32
$ws.Range("I86").Value = 1574.0625
$ws.Range("J86").Value = 2419.2222
$ws.Range("K86").Value = 1574.0625
$ws.Range("L86").Value = 2419.2222
$ws.Range("M86").Value = -451.0625
$ws.Range("N86").Value = -4665.2222
$ws.Range("H89").Value = 1878.32
$ws.Range("I89").Value = 1574.0625
$ws.Range("J89").Value = 2419.2222
$ws.Range("K89").Value = 7870.3125
$ws.Range("L89").Value = 12096.111
$ws.Range("M89").Value = -2254.3125
$ws.Range("N89").Value = -23328.111
$ws.Range("H107").Value = 2263821.5
$ws.Range("I107").Value = 2653838
$ws.Range("J107").Value = 1725.2
$ws.Range("K107").Value = 2653838
$ws.Range("L107").Value = 1725.2
$ws.Range("M107").Value = -2651918
$ws.Range("N107").Value = -5565.2
$ws.Range("H134").Value = 2404.238
$ws.Range("I134").Value = 2037
$ws.Range("K134").Value = 6111
$ws.Range("M134").Value = -3576

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4112.161
$ws.Range("I134").Value = 4828.381
$ws.Range("K134").Value = 14485.143
$ws.Range("M134").Value = -11950.143

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 999.25
$ws.Range("I93").Value = 999.25
$ws.Range("K93").Value = 2997.75
$ws.Range("M93").Value = -1125.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 47982.332
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 47982.332
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 47982.332
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -48294.332
$ws.Range("H57").Value = 34953.668
$ws.Range("I57").Value = 800
$ws.Range("J57").Value = 52030.5
$ws.Range("K57").Value = 800
$ws.Range("L57").Value = 52030.5
$ws.Range("M57").Value = 20
$ws.Range("N57").Value = -53670.5
$ws.Range("H80").Value = 100002696
$ws.Range("I80").Value = 166668880
$ws.Range("J80").Value = 3421.75
$ws.Range("K80").Value = 166668880
$ws.Range("L80").Value = 3421.75
$ws.Range("M80").Value = -166667882
$ws.Range("N80").Value = -5417.75
$ws.Range("H83").Value = 100002696
$ws.Range("I83").Value = 166668880
$ws.Range("J83").Value = 3421.75
$ws.Range("K83").Value = 833344400
$ws.Range("L83").Value = 17108.75
$ws.Range("M83").Value = -833339408
$ws.Range("N83").Value = -27092.75
$ws.Range("H132").Value = 2022.551
$ws.Range("I132").Value = 1534.8572
$ws.Range("K132").Value = 4604.571599999999
$ws.Range("M132").Value = -2074.571599999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 35000
$ws.Range("J108").Value = 35000
$ws.Range("L108").Value = 35000
$ws.Range("N108").Value = -42680
$ws.Range("H132").Value = 5025.387
$ws.Range("I132").Value = 4216
$ws.Range("K132").Value = 12648
$ws.Range("M132").Value = -10118
$ws.Range("H136").Value = 4744
$ws.Range("I136").Value = 2186.5
$ws.Range("J136").Value = 6790
$ws.Range("K136").Value = 6559.5
$ws.Range("L136").Value = 20370
$ws.Range("M136").Value = -4009.5
$ws.Range("N136").Value = -25470

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 32483.062
$ws.Range("J43").Value = 33909.285
$ws.Range("L43").Value = 33909.285
$ws.Range("N43").Value = -34207.285
$ws.Range("H51").Value = 7985
$ws.Range("I51").Value = 7985
$ws.Range("K51").Value = 7985
$ws.Range("M51").Value = -7475
$ws.Range("H52").Value = 17500
$ws.Range("I52").Value = 17500
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 17500
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -17274
$ws.Range("N52").Value = $null
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null
$ws.Range("H132").Value = 1478.8
$ws.Range("I132").Value = 1422.909
$ws.Range("K132").Value = 4268.727000000001
$ws.Range("M132").Value = -1738.727000000001
$ws.Range("H136").Value = 4933.92
$ws.Range("I136").Value = 2622.4092
$ws.Range("K136").Value = 7867.2276
$ws.Range("M136").Value = -5317.2276
